$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22/23/24: convert specific "0" text cells to real numeric 0 ---
$ws.Range("H22").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("H24").Value = 0

# --- New rows 36-52 ---
# Row 36
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = "2025-10-23"
$ws.Range("B36").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C36").Value = "CHANDERIYA `nLEAD ZINC `nSMELTER"
$ws.Range("E36").Value = "332,200 333,700 332,700 331,700 330,200"
$ws.Range("I36").NumberFormat = "@"
$ws.Range("I36").Value = "204,800"

# Row 37
$ws.Range("A37").NumberFormat = "@"
$ws.Range("A37").Value = "2025-10-23"
$ws.Range("B37").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C37").Value = "HYDRO-1 UNIT"
$ws.Range("E37").Value = "332,200 333,700 332,700 331,700 330,200"
$ws.Range("I37").NumberFormat = "@"
$ws.Range("I37").Value = "204,800"

# Row 38
$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = "2025-10-23"
$ws.Range("B38").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C38").Value = "NEW HYDRO `nSMELTER `nCHANDERIYA"
$ws.Range("E38").Value = "332,200 333,700 332,700 331,700 330,200"
$ws.Range("I38").NumberFormat = "@"
$ws.Range("I38").Value = "204,800"

# Row 39
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "2025-10-23"
$ws.Range("B39").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C39").Value = "ZINC SMELTER `nDEBRI"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0"
$ws.Range("G39").Value = "0  331,700"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "0"
$ws.Range("I39").NumberFormat = "@"
$ws.Range("I39").Value = "0"

# Row 40
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = "2025-10-23"
$ws.Range("B40").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C40").Value = "Pantnagar `nMelting&Castin `ngPlant"
$ws.Range("E40").Value = "332,200 333,700 332,700 331,700"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "0"
$ws.Range("I40").NumberFormat = "@"
$ws.Range("I40").Value = "204,800"

# Row 41
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = "2025-10-23"
$ws.Range("B41").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C41").Value = "RAJPURA DARIBA `nLEAD SMELTER"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0"
$ws.Range("F41").NumberFormat = "@"
$ws.Range("F41").Value = "0"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "0"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "0"
$ws.Range("I41").NumberFormat = "@"
$ws.Range("I41").Value = "204,800"

# Row 42
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = "2025-10-23"
$ws.Range("B42").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C42").Value = "Faridabad `nDepot"
$ws.Range("E42").Value = "334,700 336,200 330,200 334,200 332,700"
$ws.Range("I42").NumberFormat = "@"
$ws.Range("I42").Value = "207,300"

# Row 43
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = "2025-10-23"
$ws.Range("B43").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C43").Value = "Panvel Depot"
$ws.Range("E43").Value = "335,500 337,000 336,000 335,000 333,500"
$ws.Range("I43").NumberFormat = "@"
$ws.Range("I43").Value = "207,700"

# Row 44
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = "2025-10-23"
$ws.Range("B44").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C44").Value = "Pune Depot"
$ws.Range("E44").Value = "335,500 337,000 336,000 335,000 333,500"
$ws.Range("I44").NumberFormat = "@"
$ws.Range("I44").Value = "208,100"

# Row 45
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "2025-10-23"
$ws.Range("B45").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C45").Value = "Baroda Depot"
$ws.Range("E45").Value = "335,500 337,000 336,000 335,000 333,500"
$ws.Range("I45").NumberFormat = "@"
$ws.Range("I45").Value = "208,100"

# Row 46
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = "2025-10-23"
$ws.Range("B46").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C46").Value = "Raipur Depot"
$ws.Range("E46").Value = "335,500 337,000 336,000 335,000 333,500"
$ws.Range("I46").NumberFormat = "@"
$ws.Range("I46").Value = "208,100"

# Row 47
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = "2025-10-23"
$ws.Range("B47").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C47").Value = "JAMSHEDPUR `nSTOCK POINT"
$ws.Range("E47").Value = "333,200 334,700 333,700 332,700 331,200"
$ws.Range("I47").NumberFormat = "@"
$ws.Range("I47").Value = "205,800"

# Row 48
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = "2025-10-23"
$ws.Range("B48").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("D48").Value = "Kolkata Depot  333,200 334,700 333,700 332,700 331,200"
$ws.Range("I48").NumberFormat = "@"
$ws.Range("I48").Value = "205,800"

# Row 49
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = "2025-10-23"
$ws.Range("B49").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C49").Value = "Bangalore `nDepot"
$ws.Range("E49").Value = "333,200 334,700 333,700 332,700 331,200"
$ws.Range("I49").NumberFormat = "@"
$ws.Range("I49").Value = "205,800"

# Row 50
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = "2025-10-23"
$ws.Range("B50").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C50").Value = "Hyderabad `nDepot"
$ws.Range("E50").Value = "333,200 334,700 333,700 332,700 331,200"
$ws.Range("I50").NumberFormat = "@"
$ws.Range("I50").Value = "205,800"

# Row 51
$ws.Range("A51").NumberFormat = "@"
$ws.Range("A51").Value = "2025-10-23"
$ws.Range("B51").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("D51").Value = "Chennai Depot  333,200 334,700 333,700 332,700 331,200"
$ws.Range("I51").NumberFormat = "@"
$ws.Range("I51").Value = "205,800"

# Row 52
$ws.Range("A52").NumberFormat = "@"
$ws.Range("A52").Value = "2025-10-23"
$ws.Range("B52").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("C52").Value = "Sindesar `nsmelter HZAPL"
$ws.Range("E52").Value = "332,200 333,700"
$ws.Range("G52").Value = "0  331,700 330,200"
$ws.Range("I52").NumberFormat = "@"
$ws.Range("I52").Value = "204,800"

# Re-fit row heights for the newly added rows so multi-line entries don't
# leave an explicit custom row height behind.
$ws.Range("A36:I52").Rows.AutoFit()

